# Fruta / hortaliza, semanal
# Insert two new weekly observation rows into the daily-logic subset sheet
# for "Feria Lagunitas de Puerto Montt - Pepino ensalada".
#
# The sheet is a flat list of observations (row 1 = header, rows 2..203 =
# data). Two new rows are inserted:
#   - a new row at position 122 (pushing old rows 122-203 down by one)
#   - a new row at position 196, i.e. right after the (already shifted)
#     old row 194 (pushing the remaining rows down by one more)
# giving a final data range of A1:R205.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at 122 -----------------------------------------
$ws.Rows.Item(122).Insert()

$ws.Cells.Item(122, 1).Value = 4
$ws.Cells.Item(122, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(122, 3).Value = "Los Lagos"
$ws.Cells.Item(122, 4).Value = 44567
$ws.Cells.Item(122, 5).Value = 10
$ws.Cells.Item(122, 6).Value = 100112043
$ws.Cells.Item(122, 7).Value = "Pepino ensalada"
$ws.Cells.Item(122, 8).Value = "Sin especificar"
$ws.Cells.Item(122, 9).Value = "Primera"
$ws.Cells.Item(122, 10).Value = 200
$ws.Cells.Item(122, 11).Value = 10000
$ws.Cells.Item(122, 12).Value = 10000
$ws.Cells.Item(122, 13).Value = 10000
$ws.Cells.Item(122, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(122, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(122, 16).Value = 167
$ws.Cells.Item(122, 17).Value = 60
$ws.Cells.Item(122, 18).Value = "Hortaliza"

# D122 keeps the date number format used by the rest of column D
$ws.Cells.Item(122, 4).NumberFormat = $ws.Cells.Item(121, 4).NumberFormat

# --- Insert second new row at 196 (after the shifted old row 194) -------
$ws.Rows.Item(196).Insert()

$ws.Cells.Item(196, 1).Value = 4
$ws.Cells.Item(196, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(196, 3).Value = "Los Lagos"
$ws.Cells.Item(196, 4).Value = 44568
$ws.Cells.Item(196, 5).Value = 10
$ws.Cells.Item(196, 6).Value = 100112043
$ws.Cells.Item(196, 7).Value = "Pepino ensalada"
$ws.Cells.Item(196, 8).Value = "Sin especificar"
$ws.Cells.Item(196, 9).Value = "Primera"
$ws.Cells.Item(196, 10).Value = 400
$ws.Cells.Item(196, 11).Value = 12000
$ws.Cells.Item(196, 12).Value = 12000
$ws.Cells.Item(196, 13).Value = 12000
$ws.Cells.Item(196, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(196, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(196, 16).Value = 200
$ws.Cells.Item(196, 17).Value = 60
$ws.Cells.Item(196, 18).Value = "Hortaliza"

# D196 keeps the date number format used by the rest of column D
$ws.Cells.Item(196, 4).NumberFormat = $ws.Cells.Item(195, 4).NumberFormat
